$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '50.073.79'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +4.42%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.672.00'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +7.68%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '114.18'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +8.81%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '326.92'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +3.03%  '
$ws.Range("E7").Value = '  +1.95%  '
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.556'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +4.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.11'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +5.61%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.11'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.21%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0826'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +3.36%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.39'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +5.02%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.085.13'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +7.64%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.671.88'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +7.46%  '
$ws.Range("E17").Value = '  +6.48%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '49.995.33'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +4.51%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.25'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +4.62%  '
$ws.Range("E20").Value = '  +4.41%  '
$ws.Range("E21").Value = '  -2.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0963'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +3.93%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '72.62'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +2.75%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '278.02'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.91%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.60'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +4.30%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.93'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +5.18%  '
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.11'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +5.42%  '
$ws.Range("E29").Value = '  +1.80%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.52'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +5.51%  '
$ws.Range("E31").Value = '  +4.49%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '50.35'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +2.36%  '
$ws.Range("E33").Value = '  +4.75%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '19.74'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +4.21%  '
$ws.Range("E35").Value = '  +6.50%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.11'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +13.17%  '
$ws.Range("E37").Value = '  -0.04%  '
$ws.Range("E38").Value = '  +7.65%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.14'
$ws.Range("D39").ClearFormats()
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '125.27'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +2.23%  '
$ws.Range("E41").Value = '  +2.44%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '22.83'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +3.88%  '
$ws.Range("E43").Value = '  +0.52%  '
$ws.Range("E44").Value = '  +6.32%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.115.12'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +6.10%  '
$ws.Range("E46").Value = '  +5.76%  '
$ws.Range("E47").Value = '  +14.84%  '
$ws.Range("E48").Value = '  +7.38%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.13'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +2.64%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.40'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +5.30%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '59.93'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +6.99%  '
